$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "244.63"
Set-TextValue "D4" "5.414"
Set-TextValue "D6" "3.394"
Set-TextValue "D7" "0.8083"
Set-TextValue "D8" "0.9277"
Set-TextValue "D10" "0.07435"
Set-TextValue "D11" "0.03385"
Set-TextValue "D12" "0.03036"
Set-TextValue "D13" "0.09351"
Set-TextValue "D14" "3.935"
Set-TextValue "D15" "0.001592"
Set-TextValue "D16" "0.04809"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D17" "0.0005943"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.005645"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.004156"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.0009809"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.00007705"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.658"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "KuCoinToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D23" "6.468"
$ws.Range("E23").Value = "22KuCoinTokenKCS"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.186"
$ws.Range("E24").Value = "23BTSETokenBTSE"
Set-TextValue "D26" "0.1314"
Set-TextValue "D40" "0.03943"
Set-TextValue "D41" "0.006184"
Set-TextValue "D42" "0.1075"
Set-TextValue "D43" "0.002722"
Set-TextValue "D44" "0.007495"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
Set-TextValue "D45" "0.00005131"
Set-TextValue "D47" "0.0005803"
Set-TextValue "D49" "0.002251"
